$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = '34×88='
$t.Cell(1, 2).Range.Text = '37×42='
$t.Cell(1, 3).Range.Text = '38×19='
$t.Cell(1, 4).Range.Text = '87×43='
$t.Cell(1, 5).Range.Text = '32×73='
$t.Cell(2, 1).Range.Text = '42×40='
$t.Cell(2, 2).Range.Text = '19×90='
$t.Cell(2, 3).Range.Text = '65×65='
$t.Cell(2, 4).Range.Text = '25×22='
$t.Cell(2, 5).Range.Text = '81×95='
$t.Cell(3, 1).Range.Text = '55×25='
$t.Cell(3, 2).Range.Text = '36×45='
$t.Cell(3, 3).Range.Text = '74×84='
$t.Cell(3, 4).Range.Text = '30×69='
$t.Cell(3, 5).Range.Text = '24×72='
$t.Cell(4, 1).Range.Text = '14×16='
$t.Cell(4, 2).Range.Text = '19×98='
$t.Cell(4, 3).Range.Text = '47×50='
$t.Cell(4, 4).Range.Text = '63×58='
$t.Cell(4, 5).Range.Text = '41×93='
$t.Cell(5, 1).Range.Text = '88×78='
$t.Cell(5, 2).Range.Text = '93×83='
$t.Cell(5, 3).Range.Text = '36×49='
$t.Cell(5, 4).Range.Text = '18×100='
$t.Cell(5, 5).Range.Text = '62×91='
$t.Cell(6, 1).Range.Text = '93×82='
$t.Cell(6, 2).Range.Text = '98×22='
$t.Cell(6, 3).Range.Text = '71×28='
$t.Cell(6, 4).Range.Text = '11×76='
$t.Cell(6, 5).Range.Text = '14×11='
$t.Cell(7, 1).Range.Text = '53×98='
$t.Cell(7, 2).Range.Text = '20×92='
$t.Cell(7, 3).Range.Text = '78×80='
$t.Cell(7, 4).Range.Text = '53×91='
$t.Cell(7, 5).Range.Text = '68×20='
$t.Cell(8, 1).Range.Text = '100×88='
$t.Cell(8, 2).Range.Text = '39×14='
$t.Cell(8, 3).Range.Text = '90×70='
$t.Cell(8, 4).Range.Text = '26×76='
$t.Cell(8, 5).Range.Text = '98×41='
$t.Cell(9, 1).Range.Text = '70×78='
$t.Cell(9, 2).Range.Text = '98×61='
$t.Cell(9, 3).Range.Text = '100×76='
$t.Cell(9, 4).Range.Text = '33×59='
$t.Cell(9, 5).Range.Text = '53×66='
$t.Cell(10, 1).Range.Text = '89×15='
$t.Cell(10, 2).Range.Text = '84×48='
$t.Cell(10, 3).Range.Text = '14×55='
$t.Cell(10, 4).Range.Text = '28×79='
$t.Cell(10, 5).Range.Text = '61×17='
$t.Cell(11, 1).Range.Text = '55×59='
$t.Cell(11, 2).Range.Text = '49×51='
$t.Cell(11, 3).Range.Text = '67×38='
$t.Cell(11, 4).Range.Text = '65×67='
$t.Cell(11, 5).Range.Text = '58×98='
$t.Cell(12, 1).Range.Text = '28×20='
$t.Cell(12, 2).Range.Text = '76×72='
$t.Cell(12, 3).Range.Text = '40×90='
$t.Cell(12, 4).Range.Text = '21×57='
$t.Cell(12, 5).Range.Text = '100×50='
$t.Cell(13, 1).Range.Text = '58×24='
$t.Cell(13, 2).Range.Text = '94×50='
$t.Cell(13, 3).Range.Text = '33×88='
$t.Cell(13, 4).Range.Text = '61×65='
$t.Cell(13, 5).Range.Text = '92×28='
$t.Cell(14, 1).Range.Text = '80×79='
$t.Cell(14, 2).Range.Text = '33×77='
$t.Cell(14, 3).Range.Text = '69×88='
$t.Cell(14, 4).Range.Text = '13×80='
$t.Cell(14, 5).Range.Text = '43×80='
$t.Cell(15, 1).Range.Text = '83×28='
$t.Cell(15, 2).Range.Text = '20×67='
$t.Cell(15, 3).Range.Text = '67×69='
$t.Cell(15, 4).Range.Text = '49×66='
$t.Cell(15, 5).Range.Text = '66×65='
$t.Cell(16, 1).Range.Text = '77×29='
$t.Cell(16, 2).Range.Text = '18×23='
$t.Cell(16, 3).Range.Text = '57×19='
$t.Cell(16, 4).Range.Text = '99×39='
$t.Cell(16, 5).Range.Text = '51×49='
$t.Cell(17, 1).Range.Text = '61×99='
$t.Cell(17, 2).Range.Text = '17×93='
$t.Cell(17, 3).Range.Text = '27×41='
$t.Cell(17, 4).Range.Text = '83×11='
$t.Cell(17, 5).Range.Text = '88×41='
$t.Cell(18, 1).Range.Text = '38×99='
$t.Cell(18, 2).Range.Text = '85×64='
$t.Cell(18, 3).Range.Text = '67×58='
$t.Cell(18, 4).Range.Text = '22×47='
$t.Cell(18, 5).Range.Text = '24×35='
$t.Cell(19, 1).Range.Text = '18×19='
$t.Cell(19, 2).Range.Text = '47×69='
$t.Cell(19, 3).Range.Text = '78×85='
$t.Cell(19, 4).Range.Text = '68×74='
$t.Cell(19, 5).Range.Text = '22×34='
$t.Cell(20, 1).Range.Text = '98×67='
$t.Cell(20, 2).Range.Text = '72×33='
$t.Cell(20, 3).Range.Text = '11×27='
$t.Cell(20, 4).Range.Text = '59×61='
$t.Cell(20, 5).Range.Text = '83×89='
